$d = $word.ActiveDocument

# --- Step 1: split the "general conversation..." bullet into two runs ---
# Locate the boundary right after "The chatbot should allow for some"
$range = $d.Content
$found = $range.Find.Execute("The chatbot should allow for some", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target sentence start"
}
$range.Collapse(0)

# Splitting the paragraph here keeps the first half and the (still stale)
# second half as two independent runs in two different paragraphs.
$range.InsertParagraphAfter()

# Locate the paragraph holding the (still stale) remainder text by scanning,
# since paragraph indices shifted after the split above.
$bulletIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*general conversation such as greetings*") {
        $bulletIndex = $i
        break
    }
}
if ($bulletIndex -eq 0) {
    throw "Could not locate the split-off remainder paragraph"
}
$p2 = $d.Paragraphs.Item($bulletIndex)
$p2.Range.Text = " “small talk” such saying hello and asking how they are"

# Re-join the two paragraphs back into one so the split stays within the same
# bullet, but the two differently-edited runs remain distinct (no merge).
$markRange = $d.Range($range.End, $range.End + 1)
$markRange.Delete()

# --- Step 2: add the new bullet about remembering the user's name ---
$bulletIndex2 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*such saying hello and asking how they are*") {
        $bulletIndex2 = $i
        break
    }
}
if ($bulletIndex2 -eq 0) {
    throw "Could not locate the edited bullet paragraph"
}
$editedPara = $d.Paragraphs.Item($bulletIndex2)
$editedPara.Range.InsertParagraphAfter()

$newBulletIndex = $bulletIndex2 + 1
$newPara = $d.Paragraphs.Item($newBulletIndex)
$newPara.Range.Text = "Chatbot should allow the user to tell them their name and remember it in future outputs"

Write-Host "done"
